# Append the new match row (row 61) to the Azerbaijan Premier League sheet,
# mirroring the layout/styling of the last existing data row (row 60).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 61
$templateRow = 60

# --- Values -----------------------------------------------------------
$ws.Cells.Item($newRow, 1).Value = 60
$ws.Cells.Item($newRow, 2).Value = "azerbaijan"
$ws.Cells.Item($newRow, 3).Value = "premier-league"
$ws.Cells.Item($newRow, 4).Value = "2023-2024"
$ws.Cells.Item($newRow, 5).Value = 45240.66666666666
$ws.Cells.Item($newRow, 6).Value = "Araz"
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = "Kapaz"
$ws.Cells.Item($newRow, 9).Value = 1
$ws.Cells.Item($newRow, 10).Value = 1.6
$ws.Cells.Item($newRow, 11).Value = "09/11/2023 04:12"
$ws.Cells.Item($newRow, 12).Value = 1.5
$ws.Cells.Item($newRow, 13).Value = "10/11/2023 15:57"
$ws.Cells.Item($newRow, 14).Value = 3.5
$ws.Cells.Item($newRow, 15).Value = "09/11/2023 04:12"
$ws.Cells.Item($newRow, 16).Value = 3.98
$ws.Cells.Item($newRow, 17).Value = "10/11/2023 15:57"
$ws.Cells.Item($newRow, 18).Value = 4.91
$ws.Cells.Item($newRow, 19).Value = "09/11/2023 04:12"
$ws.Cells.Item($newRow, 20).Value = 6.48
$ws.Cells.Item($newRow, 21).Value = "10/11/2023 15:57"
$ws.Cells.Item($newRow, 22).Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/araz-pfk-kapaz/CSaDTmj4/"

# --- Formatting ---------------------------------------------------------
# Column A (Indice) and column E (data_partida) carry special styles
# (bold/bordered index style, datetime number format respectively) on
# every data row. Copy those formats from the previous last row (60) so
# the new row reuses the exact same style definitions.
$ws.Cells.Item($templateRow, 1).Copy() | Out-Null
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item($templateRow, 5).Copy() | Out-Null
$ws.Cells.Item($newRow, 5).PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false
